# Auto-generated edit script: apply updated LeveProfits values per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 425.5625
$ws.Range("I39").Value = 140.81818
$ws.Range("K39").Value = 422.4545400000001
$ws.Range("M39").Value = -126.4545400000001
$ws.Range("H74").Value = 3900
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 4200
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 4200
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -6072
$ws.Range("H76").Value = 3034.95
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 3099.8572
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 3099.8572
$ws.Range("M76").Value = -2685
$ws.Range("N76").Value = -3729.8572
$ws.Range("H77").Value = 3900
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 4200
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 21000
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -30360
$ws.Range("H79").Value = 3034.95
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 3099.8572
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 3099.8572
$ws.Range("M79").Value = -1908
$ws.Range("N79").Value = -5283.8572
$ws.Range("H129").Value = 743.6429000000001
$ws.Range("I129").Value = 512.3333
$ws.Range("K129").Value = 1536.9999
$ws.Range("M129").Value = 3463.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1818.6364
$ws.Range("I45").Value = 1861.2307
$ws.Range("J45").Value = 1757.1111
$ws.Range("K45").Value = 1861.2307
$ws.Range("L45").Value = 1757.1111
$ws.Range("M45").Value = -1484.2307
$ws.Range("N45").Value = -2511.1111
$ws.Range("H122").Value = 2296.6428
$ws.Range("I122").Value = 1563.5
$ws.Range("J122").Value = 2418.8333
$ws.Range("K122").Value = 4690.5
$ws.Range("L122").Value = 7256.499899999999
$ws.Range("M122").Value = -2240.5
$ws.Range("N122").Value = -12156.4999
$ws.Range("H123").Value = 19954.545
$ws.Range("J123").Value = 19954.545
$ws.Range("L123").Value = 19954.545
$ws.Range("N123").Value = -29754.545

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 99342
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 99342
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 99342
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -99998
$ws.Range("H134").Value = 3052.9219
$ws.Range("I134").Value = 2535.5
$ws.Range("J134").Value = 5847
$ws.Range("K134").Value = 7606.5
$ws.Range("L134").Value = 17541
$ws.Range("M134").Value = -5071.5
$ws.Range("N134").Value = -22611

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 27846.076
$ws.Range("J74").Value = 27846.076
$ws.Range("L74").Value = 27846.076
$ws.Range("N74").Value = -29594.076
$ws.Range("H77").Value = 27846.076
$ws.Range("J77").Value = 27846.076
$ws.Range("L77").Value = 83538.228
$ws.Range("N77").Value = -92274.228
$ws.Range("H86").Value = 9836.5
$ws.Range("I86").Value = 8019.8096
$ws.Range("J86").Value = 15286.571
$ws.Range("K86").Value = 8019.8096
$ws.Range("L86").Value = 15286.571
$ws.Range("M86").Value = -6896.8096
$ws.Range("N86").Value = -17532.571
$ws.Range("H88").Value = 25357.75
$ws.Range("J88").Value = 25357.75
$ws.Range("L88").Value = 25357.75
$ws.Range("N88").Value = -26169.75
$ws.Range("H89").Value = 9836.5
$ws.Range("I89").Value = 8019.8096
$ws.Range("J89").Value = 15286.571
$ws.Range("K89").Value = 40099.048
$ws.Range("L89").Value = 76432.855
$ws.Range("M89").Value = -34483.048
$ws.Range("N89").Value = -87664.855
$ws.Range("H91").Value = 25357.75
$ws.Range("J91").Value = 25357.75
$ws.Range("L91").Value = 25357.75
$ws.Range("N91").Value = -28165.75
$ws.Range("H107").Value = 1688.6296
$ws.Range("I107").Value = 352.93332
$ws.Range("J107").Value = 3358.25
$ws.Range("K107").Value = 352.93332
$ws.Range("L107").Value = 3358.25
$ws.Range("M107").Value = 1567.06668
$ws.Range("N107").Value = -7198.25
$ws.Range("H135").Value = 37972.727
$ws.Range("J135").Value = 37972.727
$ws.Range("L135").Value = 37972.727
$ws.Range("N135").Value = -48112.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2678.5
$ws.Range("I59").Value = 1957
$ws.Range("J59").Value = 3400
$ws.Range("K59").Value = 5871
$ws.Range("L59").Value = 10200
$ws.Range("M59").Value = -5331
$ws.Range("N59").Value = -11280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4885.6904
$ws.Range("I70").Value = 4857.6924
$ws.Range("J70").Value = 4931.1875
$ws.Range("K70").Value = 4857.6924
$ws.Range("L70").Value = 4931.1875
$ws.Range("M70").Value = -4587.6924
$ws.Range("N70").Value = -5471.1875
$ws.Range("H73").Value = 4885.6904
$ws.Range("I73").Value = 4857.6924
$ws.Range("J73").Value = 4931.1875
$ws.Range("K73").Value = 4857.6924
$ws.Range("L73").Value = 4931.1875
$ws.Range("M73").Value = -3921.6924
$ws.Range("N73").Value = -6803.1875
$ws.Range("H107").Value = 954.2778
$ws.Range("I107").Value = 873
$ws.Range("J107").Value = 1019.3
$ws.Range("K107").Value = 873
$ws.Range("L107").Value = 1019.3
$ws.Range("M107").Value = 1047
$ws.Range("N107").Value = -4859.3
$ws.Range("H113").Value = 1321.0741
$ws.Range("I113").Value = 862.06665
$ws.Range("J113").Value = 1894.8334
$ws.Range("K113").Value = 862.06665
$ws.Range("L113").Value = 1894.8334
$ws.Range("M113").Value = 1307.93335
$ws.Range("N113").Value = -6234.8334
$ws.Range("H122").Value = 73653.36
$ws.Range("I122").Value = 126555.875
$ws.Range("J122").Value = 3116.6667
$ws.Range("K122").Value = 379667.625
$ws.Range("L122").Value = 9350.000100000001
$ws.Range("M122").Value = -377217.625
$ws.Range("N122").Value = -14250.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 15387.875
$ws.Range("I68").Value = 22200
$ws.Range("J68").Value = 4034.3333
$ws.Range("K68").Value = 22200
$ws.Range("L68").Value = 4034.3333
$ws.Range("M68").Value = -21451
$ws.Range("N68").Value = -5532.3333
$ws.Range("H69").Value = 28666.666
$ws.Range("J69").Value = 28666.666
$ws.Range("L69").Value = 28666.666
$ws.Range("N69").Value = -30288.666
$ws.Range("H71").Value = 15387.875
$ws.Range("I71").Value = 22200
$ws.Range("J71").Value = 4034.3333
$ws.Range("K71").Value = 111000
$ws.Range("L71").Value = 20171.6665
$ws.Range("M71").Value = -107256
$ws.Range("N71").Value = -27659.6665
$ws.Range("H72").Value = 28666.666
$ws.Range("J72").Value = 28666.666
$ws.Range("L72").Value = 85999.99800000001
$ws.Range("N72").Value = -94111.99800000001
$ws.Range("H122").Value = 9187
$ws.Range("I122").Value = 13126
$ws.Range("J122").Value = 3935
$ws.Range("K122").Value = 39378
$ws.Range("L122").Value = 11805
$ws.Range("M122").Value = -36928
$ws.Range("N122").Value = -16705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 198.90909
$ws.Range("I113").Value = 188.5
$ws.Range("K113").Value = 565.5
$ws.Range("M113").Value = 1604.5
$ws.Range("H122").Value = 7852681.5
$ws.Range("I122").Value = 15627925
$ws.Range("J122").Value = 77438.06
$ws.Range("K122").Value = 46883775
$ws.Range("L122").Value = 232314.18
$ws.Range("M122").Value = -46881325
$ws.Range("N122").Value = -237214.18
$ws.Range("H126").Value = 47619656
$ws.Range("I126").Value = 52632156
$ws.Range("J126").Value = 895
$ws.Range("K126").Value = 157896468
$ws.Range("L126").Value = 2685
$ws.Range("M126").Value = -157893998
$ws.Range("N126").Value = -7625
$ws.Range("H132").Value = 3849879.2
$ws.Range("I132").Value = 4881537.5
$ws.Range("J132").Value = 4607.5454
$ws.Range("K132").Value = 14644612.5
$ws.Range("L132").Value = 13822.6362
$ws.Range("M132").Value = -14642082.5
$ws.Range("N132").Value = -18882.6362
$ws.Range("H136").Value = 45458330
$ws.Range("I136").Value = 83338270
$ws.Range("J136").Value = 2392
$ws.Range("K136").Value = 250014810
$ws.Range("L136").Value = 7176
$ws.Range("M136").Value = -250012260
$ws.Range("N136").Value = -12276

Write-Output "Applied 212 cell updates across 8 sheets"
